$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column G (K) for rows 2-12
$values = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 1
    7  = 0
    8  = 1
    9  = 2
    10 = 0
    11 = 3
    12 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
